$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (values 3, 5, 7, 13 with bold/bordered style) is removed,
# shifting columns B:F left into A:E.
$ws.Range("A:A").Delete()
